# Apply updated dSF (column F) values as per repull of data / mean calculation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = 1
$ws.Range("F8").Value = 4
$ws.Range("F9").Value = -6
$ws.Range("F10").Value = -1
$ws.Range("F12").Value = -1
$ws.Range("F19").Value = -6
